# Update "想去人数" (number of people interested) figures on two sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 106
$ws1.Range("F5").Value = 2711
$ws1.Range("F6").Value = 263
$ws1.Range("F7").Value = 385

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 106
$ws4.Range("F5").Value = 2711
$ws4.Range("F6").Value = 263
$ws4.Range("F9").Value = 385
